# Updated symbol list on Sat Dec 24 07:24:26 UTC 2022 with GitHub Actions
#
# This script applies the per-row price / ranking-label refresh captured in
# the commit diff. Numeric-looking values in column D are stored as TEXT in
# the workbook (same as the original inline-string cells), so each of those
# writes is preceded by forcing the cell's number format to Text ("@") -
# otherwise Excel would silently re-interpret the literal "245.83" as a
# number and drop the original text semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- simple price (column D) refreshes -------------------------------------
Set-TextValue "D2"  "245.83"
Set-TextValue "D3"  "22.05"
Set-TextValue "D4"  "5.359"
Set-TextValue "D6"  "3.392"
Set-TextValue "D7"  "6.388"
Set-TextValue "D9"  "0.9613"

# --- rows 10-18: "One" jumps from rank 18 up to rank 10, pushing --------
# --- WazirX ... CoinExToken each down by one row, each with refreshed ---
# --- price figures and renumbered rank labels in column E. -------------
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01120"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1431"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07396"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03378"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03061"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09406"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.997"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001595"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04787"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- remaining simple price (column D) refreshes ----------------------------
Set-TextValue "D19" "0.006137"
Set-TextValue "D21" "0.0009898"
Set-TextValue "D23" "3.746"
Set-TextValue "D40" "0.04207"
Set-TextValue "D41" "0.006522"
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.005807"

# --- rank-label tweak (no longer flagged "Best in 24h") ---------------------
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-TextValue "D48" "0.03289"
